$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.011893758099813892
$ws.Range("B1").Value = -0.011893758494125155
$ws.Range("A2").Value = 0.025113678923829024
$ws.Range("B2").Value = -0.02511367932214029
$ws.Range("A3").Value = -0.061676827389066269
$ws.Range("B3").Value = 0.061676826995805678
$ws.Range("A4").Value = -0.02495028971082076
$ws.Range("B4").Value = 0.02495028931612113
$ws.Range("A5").Value = 0.053820653701090709
$ws.Range("B5").Value = -0.053820654087281032
